# Updates the crypto price/volume table on Sheet1 (columns D = Price, E = Volume(1h))
# to the latest scraped values, and re-sorts a few coin pairs whose relative
# ranking flipped (rows 16/17, 41/42, 44/45, 49/50) by swapping their
# Coin/Link/Price/Volume cells.
#
# Numeric-looking Price strings (e.g. "399.11") are written with a leading
# apostrophe so Excel keeps them as text (matching the workbook's existing
# convention of storing every Price/Volume cell as a string) instead of
# auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.436.89'
$ws.Range("E2").Value = '  +9.63%  '
$ws.Range("D3").Value = '3.235.79'
$ws.Range("E3").Value = '  +4.32%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''399.11'
$ws.Range("D6").Value = '''110.60'
$ws.Range("E6").Value = '  +7.20%  '
$ws.Range("D7").Value = '''0.554'
$ws.Range("E7").Value = '  +2.47%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '''0.626'
$ws.Range("E9").Value = '  +7.58%  '
$ws.Range("D10").Value = '''39.75'
$ws.Range("E10").Value = '  +6.81%  '
$ws.Range("D11").Value = '''0.0897'
$ws.Range("E11").Value = '  +5.18%  '
$ws.Range("E12").Value = '  +2.07%  '
$ws.Range("D13").Value = '3.738.29'
$ws.Range("E13").Value = '  +4.01%  '
$ws.Range("E14").Value = '  +2.90%  '
$ws.Range("D15").Value = '''8.08'
$ws.Range("E15").Value = '  +3.37%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").Value = '''1.07'
$ws.Range("E16").Value = '  +7.45%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.237.36'
$ws.Range("E17").Value = '  +4.27%  '
$ws.Range("D18").Value = '''10.51'
$ws.Range("E18").Value = '  -6.55%  '
$ws.Range("D19").Value = '56.135.27'
$ws.Range("E19").Value = '  +8.96%  '
$ws.Range("D20").Value = '''3.36'
$ws.Range("E20").Value = '  +2.91%  '
$ws.Range("D21").Value = '''13.15'
$ws.Range("E21").Value = '  +6.35%  '
$ws.Range("D22").Value = '''0.0000101'
$ws.Range("E22").Value = '  +4.81%  '
$ws.Range("D23").Value = '''289.37'
$ws.Range("E23").Value = '  +8.84%  '
$ws.Range("D24").Value = '''74.51'
$ws.Range("E24").Value = '  +6.47%  '
$ws.Range("E25").Value = '  +4.96%  '
$ws.Range("D26").Value = '''8.22'
$ws.Range("E26").Value = '  +1.79%  '
$ws.Range("D27").Value = '''28.39'
$ws.Range("E27").Value = '  +5.26%  '
$ws.Range("E28").Value = '  +3.59%  '
$ws.Range("E29").Value = '  +3.01%  '
$ws.Range("D30").Value = '''1.00'
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("E31").Value = '  +4.55%  '
$ws.Range("D32").Value = '''11.26'
$ws.Range("D33").Value = '''0.0497'
$ws.Range("E33").Value = '  +6.13%  '
$ws.Range("D34").Value = '''37.05'
$ws.Range("E34").Value = '  +4.73%  '
$ws.Range("E35").Value = '  +1.61%  '
$ws.Range("D36").Value = '''51.04'
$ws.Range("E36").Value = '  +1.52%  '
$ws.Range("E37").Value = '  +7.57%  '
$ws.Range("D38").Value = '''0.998'
$ws.Range("E38").Value = '  -0.19%  '
$ws.Range("D39").Value = '''3.07'
$ws.Range("E39").Value = '  +21.67%  '
$ws.Range("D40").Value = '''138.01'
$ws.Range("E40").Value = '  +6.94%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '''1.93'
$ws.Range("E41").Value = '  +2.60%  '
$ws.Range("B42").Value = 'NEARProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D42").Value = '''4.04'
$ws.Range("E42").Value = '  +10.67%  '
$ws.Range("D43").Value = '''0.287'
$ws.Range("E43").Value = '  -3.37%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '''0.118'
$ws.Range("E44").Value = '  +2.17%  '
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").Value = '''16.90'
$ws.Range("E45").Value = '  +1.24%  '
$ws.Range("D46").Value = '''22.49'
$ws.Range("E46").Value = '  +0.88%  '
$ws.Range("E48").Value = '  +1.12%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").Value = '''2.08'
$ws.Range("E49").Value = '  +40.30%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.128.88'
$ws.Range("E50").Value = '  +3.38%  '
$ws.Range("D51").Value = '3.566.80'
$ws.Range("E51").Value = '  +4.45%  '
